# productDetailPage - selecting product with and without variants
#
# Adds a new "Action2" entry (E4) for the "Product Listing" test scenario:
# navigating to productDetailPage after a catalog search. Also brings the
# formatting on E3 / D4 in line with the existing "Monospace" style already
# used for the productCatalogPage entry on D2, and moves the active
# selection to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value: E4 = "productDetailPage"
$ws.Range("E4").Value = "productDetailPage"

# Normalize formatting of E3 / D4 to match the style already used by D2
# (same Monospace font) instead of the stray duplicate font/style.
$ws.Range("D2").Copy()
$ws.Range("D4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E3").PasteSpecial(-4122)  # xlPasteFormats

# Move the active selection to C3
$ws.Range("C3").Select()
